$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 98, pushing old rows 98:151 down to 99:152.
$ws.Rows("98:98").Insert()

# Populate the newly inserted row 98. Static columns match the surrounding
# records for this market/product (Agrícola del Norte S.A. de Arica - Mango),
# copied from the now-shifted row 99 (old row 98), with the new unique
# observation values for D, M, N, O, P, S.
$ws.Range("A98").Value = 1
$ws.Range("B98").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C98").Value = "Arica y Parinacota"
$ws.Range("D98").Value = 44875
$ws.Range("E98").Value = 15
$ws.Range("F98").Value = "Fruta"
$ws.Range("G98").Value = 100108
$ws.Range("H98").Value = "Tropicales y subtropicales"
$ws.Range("I98").Value = 100108002
$ws.Range("J98").Value = "Mango"
$ws.Range("K98").Value = "Sin especificar"
$ws.Range("L98").Value = "Especial"
$ws.Range("M98").Value = 400
$ws.Range("N98").Value = 8000
$ws.Range("O98").Value = 9000
$ws.Range("P98").Value = 8500
$ws.Range("Q98").Value = "$/bandeja 4 kilos"
$ws.Range("R98").Value = "Perú"
$ws.Range("S98").Value = 2125
$ws.Range("T98").Value = 4
